# chore: adapt column header formatting to respective input file names (#7)
#
# The "before"/"after" comparison columns in this AHB-Diff sheet used the
# generic "_old" / "_new" header-name suffixes. This renames them to the
# respective format-version suffixes ("_FV2304" for the "old"/left-hand
# side, "_FV2310" for the "new"/right-hand side), wraps the whole used
# range in a proper Excel Table (so the column headers are also tracked
# as table column names), and freezes the header row for easier scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.Range("A1:U57")
$headerRow = $ws.Range("A1:U1")
$lastCol = $headerRow.Columns.Count

# --- 1. Rename the header row (row 1) --------------------------------------
# Any header ending in "_old" -> "_FV2304" (the earlier format version)
# Any header ending in "_new" -> "_FV2310" (the later format version)
# Headers without either suffix (e.g. the literal "diff" column) are left
# untouched.
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Value()

    if ($header -like "*_old") {
        $cell.Value = $header -replace "_old$", "_FV2304"
    } elseif ($header -like "*_new") {
        $cell.Value = $header -replace "_new$", "_FV2310"
    }
}

# --- 2. Freeze the header row -----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Wrap the used range in an Excel Table -------------------------------
$table = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$table.Name = "Table1"

Write-Output "Renamed headers, froze top row, and added $($table.Name) over $($table.Range.Address())"
